# Apply updated dSF (column F) values as per repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 3
    "F6"  = 4
    "F7"  = 0
    "F15" = 2
    "F16" = -2
    "F23" = 9
    "F27" = -2
    "F30" = -2
    "F31" = -1
    "F39" = 3
    "F41" = -8
    "F46" = -4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
